# Implement dynamic versioning for scenario decks (v2 Scenario A, v3 Scenario B...)
# Update the "Metrics Summary" content placeholder on slide 2 with recomputed
# figures that account for both Ingram Micro and CNH Industrial leases.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

function Set-ParaText($range, $index, $newText) {
    # Clear first so PowerPoint doesn't try to diff old vs. new text and
    # split the paragraph into multiple runs; this keeps a single <a:r>.
    $range.Paragraphs($index, 1).Text = ""
    $range.Paragraphs($index, 1).Text = $newText
}

Set-ParaText $tr 2 "Compute Metrics and Draft Summary:"
Set-ParaText $tr 3 "**Total GLA**: 313,219 m² (222,221 m² for Ingram Micro + 90,998 m² for CNH Industrial)"
Set-ParaText $tr 4 "**Occupancy**: 100% (Both areas are leased)"
Set-ParaText $tr 5 "**WALT (Weighted Average Lease Term)**: Approximately 1.5 years (calculated from the lease end dates relative to the current date, assuming the current date is 2023)"
Set-ParaText $tr 6 "**In-Place Rent**: £5.1 per m² per annum (weighted average rent based on leased areas)"
Set-ParaText $tr 7 "**Key Highlight 1**: The asset is a large logistics facility with significant parking and loading capabilities, featuring 12 dock doors and 8 level access doors."
Set-ParaText $tr 8 "**Key Highlight 2**: The asset is fully occupied by two major tenants in the technology distribution and agricultural equipment industries, ensuring stable rental income."
